$d = $word.ActiveDocument

# Rebuild the document body from a clean Open XML fragment.  Using
# InsertXML on the whole-document range lets us:
#   - merge "Author :" + " Nat" into a single run and append a new
#     ", Usain Bolt" run,
#   - merge the "Blah " / "blah" / " " / "blah" / " masterpiece " runs
#     into a single run,
#   - drop the now-stale w:proofErr spell/grammar markers that were
#     sitting between the runs above (and around "intro"),
#   - append two new "conclusions" / "references" lines (each its own
#     run with a leading <w:br/>, matching the existing material /
#     discussion runs), while leaving those two runs untouched,
#   - keep the paragraph's indentation (<w:ind w:left="360"/>) and the
#     document's sectPr exactly as they were.

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r><w:t>Author : Nat</w:t></w:r>
            <w:r><w:t>, Usain Bolt</w:t></w:r>
          </w:p>
          <w:p>
            <w:r><w:t xml:space="preserve">Blah blah blah masterpiece </w:t></w:r>
          </w:p>
          <w:p>
            <w:pPr><w:ind w:left="360"/></w:pPr>
            <w:r><w:t>intro</w:t></w:r>
            <w:r><w:br/><w:t>material</w:t></w:r>
            <w:r><w:br/><w:t>discussion</w:t></w:r>
            <w:r><w:br/><w:t>conclusions</w:t></w:r>
            <w:r><w:br/><w:t>references</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$d.Content.InsertXML($xml)
